$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column letter -> new value
$changes = @{
    2  = @{ 'G' = 5.94498; 'H' = 17.83494; 'I' = 0.4679240463447598; 'J' = 0.4679240463447597;
            'M' = 0.5158793333333334; 'N' = 1.547638; 'O' = 0.0225151453405981; 'P' = 0.0225151453405981;
            'Q' = 3.06689231908; 'R' = 27.60203087172; 'S' = 0.01053537791181303; 'T' = 0.01053537791181302 }
    3  = @{ 'G' = 5.94498; 'H' = 17.83494; 'I' = 0.4679240463447598; 'J' = 0.4679240463447597;
            'O' = 0.03025903914481093; 'P' = 0.03025903914481093;
            'Q' = 4.1217239921; 'R' = 37.0955159289; 'S' = 0.01415893203514441; 'T' = 0.01415893203514441 }
    4  = @{ 'G' = 5.94498; 'H' = 17.83494; 'I' = 0.4679240463447598; 'J' = 0.4679240463447597;
            'M' = 21.70335633333333; 'N' = 65.110069; 'O' = 0.947225815514591; 'P' = 0.947225815514591;
            'Q' = 129.02601933454; 'R' = 1161.23417401086; 'S' = 0.4432297363978023; 'T' = 0.4432297363978023 }
    5  = @{ 'H' = 8.352077; 'I' = 0.219128164447035; 'J' = 0.219128164447035;
            'M' = 0.5158793333333334; 'N' = 1.547638; 'O' = 0.0225151453405981; 'P' = 0.0225151453405981;
            'Q' = 1.436221304902889; 'R' = 12.925991744126; 'S' = 0.004933702470743473; 'T' = 0.004933702470743473 }
    6  = @{ 'H' = 8.352077; 'I' = 0.219128164447035; 'J' = 0.219128164447035;
            'O' = 0.03025903914481093; 'P' = 0.03025903914481093;
            'Q' = 1.930197474999444; 'S' = 0.006630607705733398; 'T' = 0.006630607705733398 }
    7  = @{ 'H' = 8.352077; 'I' = 0.219128164447035; 'J' = 0.219128164447035;
            'M' = 21.70335633333333; 'N' = 65.110069; 'O' = 0.947225815514591; 'P' = 0.947225815514591;
            'Q' = 60.42270108481255; 'R' = 543.804309763313; 'S' = 0.2075638542705581; 'T' = 0.2075638542705581 }
    8  = @{ 'G' = 3.976005; 'H' = 11.928015; 'I' = 0.3129477892082053; 'J' = 0.3129477892082053;
            'M' = 0.5158793333333334; 'N' = 1.547638; 'O' = 0.0225151453405981; 'P' = 0.0225151453405981;
            'Q' = 2.05113880873; 'R' = 18.46024927857; 'S' = 0.007046064958041599; 'T' = 0.007046064958041599 }
    9  = @{ 'G' = 3.976005; 'H' = 11.928015; 'I' = 0.3129477892082053; 'J' = 0.3129477892082053;
            'O' = 0.03025903914481093; 'P' = 0.03025903914481093;
            'Q' = 2.756610653225; 'R' = 24.809495879025; 'S' = 0.009469499403933124; 'T' = 0.009469499403933124 }
    10 = @{ 'G' = 3.976005; 'H' = 11.928015; 'I' = 0.3129477892082053; 'J' = 0.3129477892082053;
            'M' = 21.70335633333333; 'N' = 65.110069; 'O' = 0.947225815514591; 'P' = 0.947225815514591;
            'Q' = 86.292653298115; 'R' = 776.6338796830349; 'S' = 0.2964322248462306; 'T' = 0.2964322248462306 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
